# Applies the Trade #43 close update across the workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.6       # Current Capital
$wsSummary.Range("B4").Value = -0.4         # Total P&L $
$wsSummary.Range("B5").Value = -0.19        # Total P&L %
$wsSummary.Range("B6").Value = 43           # Total Trades
$wsSummary.Range("B8").Value = 19           # Losing Trades
$wsSummary.Range("B9").Value = 34.88        # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.59999999999999   # Capital
$wsStatus.Range("D4").Value = 43                  # Trades
$wsStatus.Range("E4").Value = -0.4                # P&L $
$wsStatus.Range("F4").Value = -0.4                # P&L %
$wsStatus.Range("G4").Value = 34.88               # Win Rate %

# --- All Trades + MarketMaking sheets (Trade #43, row 44) ---
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G44").Value = 0.28
    $ws.Range("H44").Value = "CLOSED"
    $ws.Range("I44").Value = -39.1304
    $ws.Range("J44").Value = -0.18
    $ws.Range("K44").Value = 99.59999999999999
    $ws.Range("P44").Value = "early_exit"
    $ws.Range("Q44").Value = 5.1
}
